{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document starts with a single paragraph:\n//   \"Je m'entra\u00eene juste \u00e0 faire des commit donc j'\u00e9cris n'importe quoi\n//   mais faudra qu'on fasse le cahier des charges pour la semaine\n//   prochaine.<br>T'fa\u00e7on le sujet est incompr\u00e9hensible. Je vais aussi le\n//   rajouter au d\u00e9p\u00f4t git \u00e7a sera rigolo !\" followed by the _GoBack\n// bookmark. We turn it into four paragraphs of new brainstorming notes,\n// keeping the trailing bookmark on the final paragraph.\nconst lastParagraph = paragraphs.items[0];\n\n// Insert the three new, simple paragraphs ahead of the existing one.\nlastParagraph.insertParagraph(\"Bon, il nous faut un blog\", Word.InsertLocation.before);\nlastParagraph.insertParagraph(\"Donc : Faut avoir des articles, avec des commentaires\", Word.InsertLocation.before);\nconst thirdParagraph = lastParagraph.insertParagraph(\n  \"Dans les articles faut pouvoir afficher des images, des liens ?\",\n  Word.InsertLocation.before\n);\nawait context.sync();\n\n// That third paragraph ends with a lone line break.\nthirdParagraph.insertBreak(Word.BreakType.line, Word.InsertLocation.end);\nawait context.sync();\n\n// Rewrite the original (last) paragraph's OOXML directly so the line\n// breaks land in their own runs while the _GoBack bookmark stays at the\n// very end of the paragraph.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>Faut pouvoir voir si un article est plus r\\u00e9cent qu\\u2019un autre ou nom.</w:t></w:r>\n            <w:r><w:br/></w:r>\n            <w:r><w:br/><w:t>Pour les commentaires, des utilisateurs ?</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n          <w:sectPr/>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nlastParagraph.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document starts as a single paragraph:\n#   \"Je m'entra\u00eene juste \u00e0 faire des commit donc j'\u00e9cris n'importe quoi\n#   mais faudra qu'on fasse le cahier des charges pour la semaine\n#   prochaine.<br>T'fa\u00e7on le sujet est incompr\u00e9hensible. Je vais aussi le\n#   rajouter au d\u00e9p\u00f4t git \u00e7a sera rigolo !\" followed by the _GoBack\n# bookmark. Turn it into four paragraphs of new brainstorming notes,\n# keeping the trailing bookmark attached to the final paragraph.\n\n$p1 = $d.Paragraphs(1)\n\n# Insert three blank paragraphs ahead of the existing one; this keeps\n# the original paragraph (and its _GoBack bookmark) as paragraph 4.\n$p1.Range.InsertParagraphBefore()\n$p1.Range.InsertParagraphBefore()\n$p1.Range.InsertParagraphBefore()\n\n# Paragraphs(1..3).Range.Text replaces just the leading run of text,\n# which for these brand new, single-run paragraphs gives each one the\n# exact wording we want.\n$d.Paragraphs(1).Range.Text = \"Bon, il nous faut un blog\"\n$d.Paragraphs(2).Range.Text = \"Donc : Faut avoir des articles, avec des commentaires\"\n$d.Paragraphs(3).Range.Text = \"Dans les articles faut pouvoir afficher des images, des liens ?\"\n\n# Paragraph 3 ends with a lone line break.\n$p3 = $d.Paragraphs(3)\n$endOfP3 = $p3.Range\n$endOfP3.Collapse(0)\n$endOfP3.InsertBreak(6)\n\n# Paragraph 4 is still the original paragraph, carrying the _GoBack\n# bookmark. Replace its whole range with the exact target OOXML so the\n# line breaks land in their own runs and the bookmark stays at the end.\n$r4 = $d.Paragraphs(4).Range\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>Faut pouvoir voir si un article est plus r\u00e9cent qu\u2019un autre ou nom.</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t>Pour les commentaires, des utilisateurs ?</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$r4.InsertXML($ooxml)\n"}
